# PEARL_MOOS_testing_notes.xlsx - add 2021-03-25 testing log entries (rows 45-52)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Add-LogRow {
    param($Row, $DateSerial, $TimeFraction, $Mission, $MissionStyleSrc, $Notes, $NotesStyleSrc, $RowHeight, $DateStyleSrc, $TimeStyleSrc)

    # --- Date (column A) ---
    $ws.Range("A$Row").Value = $DateSerial
    $ws.Range($DateStyleSrc).Copy()
    $ws.Range("A$Row").PasteSpecial($xlPasteFormats)

    # --- Approximate Time (column B) ---
    $ws.Range("B$Row").Value = $TimeFraction
    $ws.Range($TimeStyleSrc).Copy()
    $ws.Range("B$Row").PasteSpecial($xlPasteFormats)

    # --- Mission (column C) ---
    $ws.Range("C$Row").Value = $Mission
    $ws.Range($MissionStyleSrc).Copy()
    $ws.Range("C$Row").PasteSpecial($xlPasteFormats)

    # --- Gain columns E:J (fixed values common to all new rows) ---
    $ws.Range("E$Row").Value = 0.4
    $ws.Range("F$Row").Value = 0.04
    $ws.Range("G$Row").Value = 4
    $ws.Range("H$Row").Value = 1
    $ws.Range("I$Row").Value = 0
    $ws.Range("J$Row").Value = 0.5
    $ws.Range("E42:J42").Copy()
    $ejRange = "E$Row" + ":J$Row"
    $ws.Range($ejRange).PasteSpecial($xlPasteFormats)

    # --- Notes (column K) ---
    $ws.Range("K$Row").Value = $Notes
    $ws.Range($NotesStyleSrc).Copy()
    $ws.Range("K$Row").PasteSpecial($xlPasteFormats)

    # --- Row height ---
    $ws.Rows($Row).RowHeight = $RowHeight

    Write-Output "row $Row written"
}

$notes48 = "- virtually no wind or waves now, water quite still, seems to be following path exceptionally well!
- my anchor broke loose??"

$notes50 = "- redo for statistics
- got caught on something or went stale? Hmm"

Add-LogRow 45 44280 0.56874999999999998 "Waypoint - Figure 8" "C42" "- PEARL bumped into edge of canoe on the way to first way point (grazed). Will it show up in log?" "K40" 34 "A42" "B44"
Add-LogRow 46 44280 0.57777777777777783 "Waypoint - Figure 8" "C42" "- did well, some overshoot/undershoot on turns but overall OK" "K40" 17 "A42" "B44"
Add-LogRow 47 44280 0.58680555555555558 "Waypoint - Figure 8" "C42" "Oops no good, kept starting wrong mission" "K2" 17 "A42" "B44"
Add-LogRow 48 44280 0.58819444444444446 "Waypoint - Figure 8" "C42" $notes48 "K40" 51 "A42" "B44"
Add-LogRow 49 44280 0.59583333333333333 "Waypoint - Star" "C42" "- excellent!" "K40" 17 "A42" "B44"
# Row 51's shared string ("- redo for statistics") is interned before row 50's
# ("- redo for statistics\n- got caught...") so the shared-string table order
# matches the source workbook (short string = earlier unique index).
Add-LogRow 51 44280 0.6069444444444444 "Waypoint - Star" "C42" "- redo for statistics" "K40" 17 "A42" "B44"
Add-LogRow 50 44280 0.60555555555555551 "Waypoint - Star" "C42" $notes50 "K40" 34 "A42" "B44"
Add-LogRow 52 44280 0.61736111111111114 "Simple" "C42" "- to send her home cuz lazy lol" "K40" 17 "A40" "B44"

# --- Selection / scroll position to match the saved view ---
$win = $excel.ActiveWindow
$ws.Range("B53").Select()
$win.ScrollRow = 37
$win.ScrollColumn = 1

Write-Output "Added rows 45-52 for 2021-03-25 testing log"
